$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    # Force a numeric-looking string to be stored as text (matching the
    # original inline-string cell type) without altering the cell style:
    # build it as a formula returning the literal string, then collapse
    # the formula down to its static value via copy / paste-special.
    $escaped = $text -replace '"', '""'
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null
}

# Row 2
Set-TextValue $ws.Range("D2") "43.846.98"
$ws.Range("E2").Value = "  -0.58%  "

# Row 3
Set-TextValue $ws.Range("D3") "2.365.21"
$ws.Range("E3").Value = "  +0.43%  "

# Row 4
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("E5").Value = "  -1.35%  "

# Row 6
Set-TextValue $ws.Range("D6") "239.67"
$ws.Range("E6").Value = "  +0.19%  "

# Row 7
$ws.Range("E7").Value = "  -0.26%  "

# Row 8
$ws.Range("E8").Value = "  +0.06%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.600"
$ws.Range("E9").Value = "  +0.39%  "

# Row 10
$ws.Range("E10").Value = "  +1.58%  "

# Row 11
Set-TextValue $ws.Range("D11") "59.80"
$ws.Range("E11").Value = "  +4.34%  "

# Row 12
Set-TextValue $ws.Range("D12") "36.80"
$ws.Range("E12").Value = "  +13.56%  "

# Row 13
$ws.Range("E13").Value = "  +0.44%  "

# Row 14
Set-TextValue $ws.Range("D14") "7.28"
$ws.Range("E14").Value = "  -0.08%  "

# Row 15
Set-TextValue $ws.Range("D15") "16.36"
$ws.Range("E15").Value = "  -1.36%  "

# Row 16
Set-TextValue $ws.Range("D16") "0.929"
$ws.Range("E16").Value = "  +3.21%  "

# Row 17
Set-TextValue $ws.Range("D17") "2.372.19"
$ws.Range("E17").Value = "  +0.70%  "

# Row 18
Set-TextValue $ws.Range("D18") "43.855.90"
$ws.Range("E18").Value = "  -0.20%  "

# Row 19
$ws.Range("E19").Value = "  +1.39%  "

# Row 20
$ws.Range("E20").Value = "  -1.87%  "

# Row 21
Set-TextValue $ws.Range("D21") "77.46"
$ws.Range("E21").Value = "  +0.70%  "

# Row 22
Set-TextValue $ws.Range("D22") "254.03"
$ws.Range("E22").Value = "  -2.06%  "

# Row 23
$ws.Range("E23").Value = "  -0.12%  "

# Row 24
$ws.Range("E24").Value = "  +3.78%  "

# Row 25
Set-TextValue $ws.Range("D25") "1.89"
$ws.Range("E25").Value = "  -3.97%  "

# Row 26
$ws.Range("E26").Value = "  +0.15%  "

# Row 27
$ws.Range("E27").Value = "  -1.53%  "

# Row 28
$ws.Range("E28").Value = "  +0.86%  "

# Row 29
$ws.Range("E29").Value = "  -1.46%  "

# Row 30
Set-TextValue $ws.Range("D30") "175.60"
$ws.Range("E30").Value = "  -0.24%  "

# Row 31
$ws.Range("E31").Value = "  +0.64%  "

# Row 32
$ws.Range("E32").Value = "  -1.27%  "

# Row 33
Set-TextValue $ws.Range("D33") "0.0759"
$ws.Range("E33").Value = "  -0.25%  "

# Row 34
Set-TextValue $ws.Range("D34") "5.42"
$ws.Range("E34").Value = "  -1.90%  "

# Row 35
$ws.Range("E35").Value = "  -2.84%  "

# Row 36
Set-TextValue $ws.Range("D36") "3.80"
$ws.Range("E36").Value = "  +1.05%  "

# Row 37
$ws.Range("E37").Value = "  +4.90%  "

# Row 38
$ws.Range("E38").Value = "  +2.15%  "

# Row 39
$ws.Range("E39").Value = "  +1.32%  "

# Row 40
Set-TextValue $ws.Range("D40") "5.56"
$ws.Range("E40").Value = "  +18.72%  "

# Row 41
Set-TextValue $ws.Range("D41") "20.49"
$ws.Range("E41").Value = "  +8.12%  "

# Row 42
Set-TextValue $ws.Range("D42") "65.63"
$ws.Range("E42").Value = "  +12.73%  "

# Row 43
$ws.Range("E43").Value = "  -3.35%  "

# Row 44
$ws.Range("E44").Value = "  -0.96%  "

# Row 45
Set-TextValue $ws.Range("D45") "9.08"
$ws.Range("E45").Value = "  +1.37%  "

# Row 46
Set-TextValue $ws.Range("D46") "2.56"
$ws.Range("E46").Value = "  +2.01%  "

# Row 47
$ws.Range("E47").Value = "  +0.54%  "

# Row 48
$ws.Range("E48").Value = "  +0.07%  "

# Row 49
$ws.Range("E49").Value = "  -1.27%  "

# Row 50
Set-TextValue $ws.Range("D50") "98.41"
$ws.Range("E50").Value = "  -1.79%  "

# Row 51
$ws.Range("B51").Value = "HuobiToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws.Range("D51") "2.88"
$ws.Range("E51").Value = "  +2.14%  "

$excel.CutCopyMode = 0